# Set cell E4 to the value 4 and select it, matching the authored edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 4
$ws.Range("E4").Select() | Out-Null
